$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM recomputation expanded the sending/target cluster
# combinations from 2x3 to the full 3x3 permutation of ECs/FAPs/MuSCs, and
# refreshed every derived metric (columns E:T) for each row accordingly.

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dcn"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.006697333333333
$ws.Range("H2").Value = 3.020092
$ws.Range("I2").Value = 0.0001985651645046208
$ws.Range("J2").Value = 0.0001985651645046208
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.493155
$ws.Range("N2").Value = 1.479465
$ws.Range("O2").Value = 0.2262533155038342
$ws.Range("P2").Value = 0.2262533155038342
$ws.Range("Q2").Value = 0.49645782342
$ws.Range("R2").Value = 4.46812041078
$ws.Range("S2").Value = [double]"4.492602681273471E-05"
$ws.Range("T2").Value = [double]"4.492602681273471E-05"

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dcn"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.006697333333333
$ws.Range("H3").Value = 3.020092
$ws.Range("I3").Value = 0.0001985651645046208
$ws.Range("J3").Value = 0.0001985651645046208
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.648742666666666
$ws.Range("N3").Value = 4.946228
$ws.Range("O3").Value = 0.7564224123165462
$ws.Range("P3").Value = 0.7564224123165462
$ws.Range("Q3").Value = 1.659784845886222
$ws.Range("R3").Value = 14.938063612976
$ws.Range("S3").Value = 0.0001501991407366171
$ws.Range("T3").Value = 0.0001501991407366171

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dcn"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.006697333333333
$ws.Range("H4").Value = 3.020092
$ws.Range("I4").Value = 0.0001985651645046208
$ws.Range("J4").Value = 0.0001985651645046208
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.037761
$ws.Range("N4").Value = 0.113283
$ws.Range("O4").Value = 0.01732427217961956
$ws.Range("P4").Value = 0.01732427217961956
$ws.Range("Q4").Value = 0.038013898004
$ws.Range("R4").Value = 0.342125082036
$ws.Range("S4").Value = [double]"3.439996955268984E-06"
$ws.Range("T4").Value = [double]"3.439996955268984E-06"

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dcn"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4971.754394666666
$ws.Range("H5").Value = 14915.263184
$ws.Range("I5").Value = 0.9806494927176636
$ws.Range("J5").Value = 0.9806494927176637
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.493155
$ws.Range("N5").Value = 1.479465
$ws.Range("O5").Value = 0.2262533155038342
$ws.Range("P5").Value = 0.2262533155038342
$ws.Range("Q5").Value = 2451.84553850184
$ws.Range("R5").Value = 22066.60984651656
$ws.Range("S5").Value = 0.2218751990745245
$ws.Range("T5").Value = 0.2218751990745245

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dcn"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4971.754394666666
$ws.Range("H6").Value = 14915.263184
$ws.Range("I6").Value = 0.9806494927176636
$ws.Range("J6").Value = 0.9806494927176637
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.648742666666666
$ws.Range("N6").Value = 4.946228
$ws.Range("O6").Value = 0.7564224123165462
$ws.Range("P6").Value = 0.7564224123165462
$ws.Range("Q6").Value = 8197.143598674438
$ws.Range("R6").Value = 73774.29238806994
$ws.Range("S6").Value = 0.7417852549184923
$ws.Range("T6").Value = 0.7417852549184925

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dcn"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4971.754394666666
$ws.Range("H7").Value = 14915.263184
$ws.Range("I7").Value = 0.9806494927176636
$ws.Range("J7").Value = 0.9806494927176637
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.037761
$ws.Range("N7").Value = 0.113283
$ws.Range("O7").Value = 0.01732427217961956
$ws.Range("P7").Value = 0.01732427217961956
$ws.Range("Q7").Value = 187.738417697008
$ws.Range("R7").Value = 1689.645759273072
$ws.Range("S7").Value = 0.01698903872464665
$ws.Range("T7").Value = 0.01698903872464666

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Dcn"
$ws.Range("C8").Value = "Tlr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 97.097641
$ws.Range("H8").Value = 291.292923
$ws.Range("I8").Value = 0.01915194211783179
$ws.Range("J8").Value = 0.01915194211783179
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.493155
$ws.Range("N8").Value = 1.479465
$ws.Range("O8").Value = 0.2262533155038342
$ws.Range("P8").Value = 0.2262533155038342
$ws.Range("Q8").Value = 47.884187147355
$ws.Range("R8").Value = 430.957684326195
$ws.Range("S8").Value = 0.004333190402496967
$ws.Range("T8").Value = 0.004333190402496967

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Dcn"
$ws.Range("C9").Value = "Tlr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 97.097641
$ws.Range("H9").Value = 291.292923
$ws.Range("I9").Value = 0.01915194211783179
$ws.Range("J9").Value = 0.01915194211783179
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.648742666666666
$ws.Range("N9").Value = 4.946228
$ws.Range("O9").Value = 0.7564224123165462
$ws.Range("P9").Value = 0.7564224123165462
$ws.Range("Q9").Value = 160.0890235493826
$ws.Range("R9").Value = 1440.801211944444
$ws.Range("S9").Value = 0.01448695825731718
$ws.Range("T9").Value = 0.01448695825731718

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Dcn"
$ws.Range("C10").Value = "Tlr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 97.097641
$ws.Range("H10").Value = 291.292923
$ws.Range("I10").Value = 0.01915194211783179
$ws.Range("J10").Value = 0.01915194211783179
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.037761
$ws.Range("N10").Value = 0.113283
$ws.Range("O10").Value = 0.01732427217961956
$ws.Range("P10").Value = 0.01732427217961956
$ws.Range("Q10").Value = 3.666504021800999
$ws.Range("R10").Value = 32.998536196209
$ws.Range("S10").Value = 0.0003317934580176374
$ws.Range("T10").Value = 0.0003317934580176374
